# Apply scraped leve-profit updates across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H22").Value = 10000
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 30000
$ws.Range("N22").Value = -30344

$ws.Range("H38").Value = 299.5
$ws.Range("I38").Value = 68.69231
$ws.Range("J38").Value = 3300
$ws.Range("K38").Value = 206.07693
$ws.Range("L38").Value = 9900
$ws.Range("M38").Value = 165.92307
$ws.Range("N38").Value = -10644

$ws.Range("H58").Value = 2029.9375
$ws.Range("I58").Value = 179.16667
$ws.Range("J58").Value = 3140.4
$ws.Range("K58").Value = 537.50001
$ws.Range("L58").Value = 9421.2
$ws.Range("M58").Value = -387.50001
$ws.Range("N58").Value = -9721.2

$ws.Range("H64").Value = 4197.6665
$ws.Range("J64").Value = 4397.206
$ws.Range("L64").Value = 4397.206
$ws.Range("N64").Value = -4893.206

$ws.Range("H67").Value = 4197.6665
$ws.Range("J67").Value = 4397.206
$ws.Range("L67").Value = 4397.206
$ws.Range("N67").Value = -6113.206

$ws.Range("H76").Value = 4100.5557
$ws.Range("I76").Value = 3992.5
$ws.Range("K76").Value = 3992.5
$ws.Range("M76").Value = -3677.5

$ws.Range("H79").Value = 4100.5557
$ws.Range("I79").Value = 3992.5
$ws.Range("K79").Value = 3992.5
$ws.Range("M79").Value = -2900.5

$ws.Range("H104").Value = 290
$ws.Range("I104").Value = 290
$ws.Range("K104").Value = 870
$ws.Range("M104").Value = 877

$ws.Range("H112").Value = 2070.7837
$ws.Range("J112").Value = 2212.6667
$ws.Range("L112").Value = 6638.000100000001
$ws.Range("N112").Value = -8854.000100000001

$ws.Range("H129").Value = 960.8
$ws.Range("I129").Value = 469.35715
$ws.Range("J129").Value = 1182.742
$ws.Range("K129").Value = 1408.07145
$ws.Range("L129").Value = 3548.226
$ws.Range("M129").Value = 3591.92855
$ws.Range("N129").Value = -13548.226

$ws.Range("H138").Value = 2945804
$ws.Range("I138").Value = 6899102.5
$ws.Range("J138").Value = 6171.769
$ws.Range("K138").Value = 20697307.5
$ws.Range("L138").Value = 18515.307
$ws.Range("M138").Value = -20692167.5
$ws.Range("N138").Value = -28795.307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15830.482
$ws.Range("I32").Value = 14736.306
$ws.Range("J32").Value = 23489.715
$ws.Range("K32").Value = 14736.306
$ws.Range("L32").Value = 23489.715
$ws.Range("M32").Value = -14449.306
$ws.Range("N32").Value = -24063.715

$ws.Range("H88").Value = 2717.8333
$ws.Range("I88").Value = 1850
$ws.Range("J88").Value = 3151.75
$ws.Range("K88").Value = 1850
$ws.Range("L88").Value = 3151.75
$ws.Range("M88").Value = -1444
$ws.Range("N88").Value = -3963.75

$ws.Range("H91").Value = 2717.8333
$ws.Range("I91").Value = 1850
$ws.Range("J91").Value = 3151.75
$ws.Range("K91").Value = 1850
$ws.Range("L91").Value = 3151.75
$ws.Range("M91").Value = -446
$ws.Range("N91").Value = -5959.75

$ws.Range("H110").Value = 1009.5333
$ws.Range("I110").Value = 874.3
$ws.Range("J110").Value = 1280
$ws.Range("K110").Value = 874.3
$ws.Range("L110").Value = 1280
$ws.Range("M110").Value = 1170.7
$ws.Range("N110").Value = -5370

$ws.Range("H132").Value = 489507.84
$ws.Range("I132").Value = 646468.3
$ws.Range("J132").Value = 2930.4
$ws.Range("K132").Value = 1939404.9
$ws.Range("L132").Value = 8791.2
$ws.Range("M132").Value = -1936874.9
$ws.Range("N132").Value = -13851.2

$ws.Range("H141").Value = 58900
$ws.Range("J141").Value = 58900
$ws.Range("L141").Value = 58900
$ws.Range("N141").Value = -69260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 12897.5
$ws.Range("I22").Value = 17098.334
$ws.Range("K22").Value = 17098.334
$ws.Range("M22").Value = -16925.334

$ws.Range("H94").Value = 1386.1578
$ws.Range("I94").Value = 1059.7858
$ws.Range("J94").Value = 2300
$ws.Range("K94").Value = 1059.7858
$ws.Range("L94").Value = 2300
$ws.Range("M94").Value = -608.7858000000001
$ws.Range("N94").Value = -3202

$ws.Range("H105").Value = 3629.5334
$ws.Range("I105").Value = 3174.5
$ws.Range("K105").Value = 3174.5
$ws.Range("M105").Value = -1427.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1389.2413
$ws.Range("I58").Value = 1535.2222
$ws.Range("J58").Value = 1150.3636
$ws.Range("K58").Value = 1535.2222
$ws.Range("L58").Value = 1150.3636
$ws.Range("M58").Value = -1332.2222
$ws.Range("N58").Value = -1556.3636

$ws.Range("H62").Value = 57941.11
$ws.Range("I62").Value = 85686.836
$ws.Range("K62").Value = 85686.836
$ws.Range("M62").Value = -85062.836

$ws.Range("H65").Value = 57941.11
$ws.Range("I65").Value = 85686.836
$ws.Range("K65").Value = 428434.18
$ws.Range("M65").Value = -425314.18

$ws.Range("H132").Value = 2993.7144
$ws.Range("I132").Value = 2690.7693
$ws.Range("J132").Value = 3486
$ws.Range("K132").Value = 8072.3079
$ws.Range("L132").Value = 10458
$ws.Range("M132").Value = -5542.3079
$ws.Range("N132").Value = -15518

$ws.Range("H134").Value = 1115.2
$ws.Range("I134").Value = 953.0278
$ws.Range("J134").Value = 2574.75
$ws.Range("K134").Value = 2859.0834
$ws.Range("L134").Value = 7724.25
$ws.Range("M134").Value = -324.0834
$ws.Range("N134").Value = -12794.25

$ws.Range("H136").Value = 1389.2413
$ws.Range("I136").Value = 1535.2222
$ws.Range("J136").Value = 1150.3636
$ws.Range("K136").Value = 4605.6666
$ws.Range("L136").Value = 3451.0908
$ws.Range("M136").Value = -2055.6666
$ws.Range("N136").Value = -8551.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4231.567
$ws.Range("I113").Value = 594.3333
$ws.Range("J113").Value = 5790.381
$ws.Range("K113").Value = 1782.9999
$ws.Range("L113").Value = 17371.143
$ws.Range("M113").Value = 387.0001
$ws.Range("N113").Value = -21711.143

$ws.Range("H122").Value = 3548.543
$ws.Range("I122").Value = 441.1
$ws.Range("J122").Value = 7691.8
$ws.Range("K122").Value = 3969.9
$ws.Range("L122").Value = 69226.2
$ws.Range("M122").Value = -1519.9
$ws.Range("N122").Value = -74126.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5776.6553
$ws.Range("I70").Value = 4894.2144
$ws.Range("J70").Value = 6600.2666
$ws.Range("K70").Value = 4894.2144
$ws.Range("L70").Value = 6600.2666
$ws.Range("M70").Value = -4624.2144
$ws.Range("N70").Value = -7140.2666

$ws.Range("H73").Value = 5776.6553
$ws.Range("I73").Value = 4894.2144
$ws.Range("J73").Value = 6600.2666
$ws.Range("K73").Value = 4894.2144
$ws.Range("L73").Value = 6600.2666
$ws.Range("M73").Value = -3958.2144
$ws.Range("N73").Value = -8472.266599999999

$ws.Range("H80").Value = 3135.9092
$ws.Range("I80").Value = 2857.8572
$ws.Range("J80").Value = 3622.5
$ws.Range("K80").Value = 2857.8572
$ws.Range("L80").Value = 3622.5
$ws.Range("M80").Value = -1859.8572
$ws.Range("N80").Value = -5618.5

$ws.Range("H83").Value = 3135.9092
$ws.Range("I83").Value = 2857.8572
$ws.Range("J83").Value = 3622.5
$ws.Range("K83").Value = 14289.286
$ws.Range("L83").Value = 18112.5
$ws.Range("M83").Value = -9297.286
$ws.Range("N83").Value = -28096.5

$ws.Range("H126").Value = 2859.7778
$ws.Range("I126").Value = 1635.2
$ws.Range("J126").Value = 3330.7693
$ws.Range("K126").Value = 4905.6
$ws.Range("L126").Value = 9992.3079
$ws.Range("M126").Value = -2435.6
$ws.Range("N126").Value = -14932.3079

$ws.Range("H132").Value = 2374.0952
$ws.Range("I132").Value = 1351.1666
$ws.Range("J132").Value = 3738
$ws.Range("K132").Value = 4053.4998
$ws.Range("L132").Value = 11214
$ws.Range("M132").Value = -1523.4998
$ws.Range("N132").Value = -16274

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 919.8333
$ws.Range("I16").Value = 1084
$ws.Range("K16").Value = 1084
$ws.Range("M16").Value = -914

$ws.Range("H140").Value = 76202.5
$ws.Range("J140").Value = 76202.5
$ws.Range("L140").Value = 76202.5
$ws.Range("N140").Value = -86562.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 27817.166
$ws.Range("J103").Value = 27817.166
$ws.Range("L103").Value = 27817.166
$ws.Range("N103").Value = -30161.166

$ws.Range("H107").Value = 735.4545
$ws.Range("I107").Value = 716.6667
$ws.Range("J107").Value = 758
$ws.Range("K107").Value = 2150.0001
$ws.Range("L107").Value = 2274
$ws.Range("M107").Value = -230.0001000000002
$ws.Range("N107").Value = -6114

$ws.Range("H136").Value = 1767.6305
$ws.Range("I136").Value = 1628
$ws.Range("J136").Value = 2341.6667
$ws.Range("K136").Value = 4884
$ws.Range("L136").Value = 7025.000100000001
$ws.Range("M136").Value = -2334
$ws.Range("N136").Value = -12125.0001
